# The menu's numbering column (A/B/C) held small local sequence numbers
# (1, 2, 3, ...). They are replaced with global catalog id numbers stored
# as text (e.g. "12", "17", "37", ...).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$idMap = @{
    "A1"  = "12"
    "B2"  = "17"
    "C3"  = "37"
    "C4"  = "38"
    "C5"  = "39"
    "B6"  = "18"
    "C7"  = "40"
    "C8"  = "41"
    "C9"  = "42"
    "A10" = "13"
    "B11" = "19"
    "C12" = "43"
    "C13" = "44"
    "C14" = "45"
    "B15" = "20"
    "C16" = "46"
    "C17" = "47"
    "C18" = "48"
}

foreach ($addr in $idMap.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $idMap[$addr]
}
